$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("G2").Value = 0.7454756666666666
$ws.Range("H2").Value = 2.236427
$ws.Range("I2").Value = 0.004237455828433692
$ws.Range("J2").Value = 0.004251944035061194
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 30.34606566666666
$ws.Range("N2").Value = 91.038197
$ws.Range("O2").Value = 0.7437342022026642
$ws.Range("P2").Value = 0.744350442811843
$ws.Range("Q2").Value = 22.62225353356877
$ws.Range("R2").Value = 203.600281802119
$ws.Range("S2").Value = 0.003151540829929161
$ws.Range("T2").Value = 0.003164936425308974
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 0.7454756666666666
$ws.Range("H3").Value = 2.236427
$ws.Range("I3").Value = 0.004237455828433692
$ws.Range("J3").Value = 0.004251944035061194
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.1013395
$ws.Range("N3").Value = 0.202679
$ws.Range("O3").Value = 0.002483671294065179
$ws.Range("P3").Value = 0.001657152803659562
$ws.Range("Q3").Value = 0.07554613132216667
$ws.Range("R3").Value = 0.453276787933
$ws.Range("S3").Value = [double]"1.052444740094994E-05"
$ws.Range("T3").Value = [double]"7.046120978705209E-06"
$ws.Range("G4").Value = 0.7454756666666666
$ws.Range("H4").Value = 2.236427
$ws.Range("I4").Value = 0.004237455828433692
$ws.Range("J4").Value = 0.004251944035061194
$ws.Range("M4").Value = 10.35489433333333
$ws.Range("N4").Value = 31.064683
$ws.Range("O4").Value = 0.2537821265032705
$ws.Range("P4").Value = 0.2539924043844974
$ws.Range("Q4").Value = 7.719321756404556
$ws.Range("R4").Value = 69.473895807641
$ws.Range("S4").Value = 0.00107539055110358
$ws.Range("T4").Value = 0.001079961488773514
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("G5").Value = 102.35201
$ws.Range("H5").Value = 307.05603
$ws.Range("I5").Value = 0.5817924591230612
$ws.Range("J5").Value = 0.583781654929077
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 30.34606566666666
$ws.Range("N5").Value = 91.038197
$ws.Range("O5").Value = 0.7437342022026642
$ws.Range("P5").Value = 0.744350442811843
$ws.Range("Q5").Value = 3105.980816575323
$ws.Range("R5").Value = 27953.82734917791
$ws.Range("S5").Value = 0.4326989504334161
$ws.Range("T5").Value = 0.434538133351889
$ws.Range("D6").Value = "MuSCs"
$ws.Range("G6").Value = 102.35201
$ws.Range("H6").Value = 307.05603
$ws.Range("I6").Value = 0.5817924591230612
$ws.Range("J6").Value = 0.583781654929077
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.1013395
$ws.Range("N6").Value = 0.202679
$ws.Range("O6").Value = 0.002483671294065179
$ws.Range("P6").Value = 0.001657152803659562
$ws.Range("Q6").Value = 10.372301517395
$ws.Range("R6").Value = 62.23380910436999
$ws.Range("S6").Value = 0.001444981229827536
$ws.Range("T6").Value = 0.0009674154061907389
$ws.Range("G7").Value = 102.35201
$ws.Range("H7").Value = 307.05603
$ws.Range("I7").Value = 0.5817924591230612
$ws.Range("J7").Value = 0.583781654929077
$ws.Range("M7").Value = 10.35489433333333
$ws.Range("N7").Value = 31.064683
$ws.Range("O7").Value = 0.2537821265032705
$ws.Range("P7").Value = 0.2539924043844974
$ws.Range("Q7").Value = 1059.844248354277
$ws.Range("R7").Value = 9538.59823518849
$ws.Range("S7").Value = 0.1476485274598176
$ws.Range("T7").Value = 0.1482761061709972
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("G8").Value = 1.79836
$ws.Range("H8").Value = 3.59672
$ws.Range("I8").Value = 0.01022229350247785
$ws.Range("J8").Value = 0.006838162904394061
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 30.34606566666666
$ws.Range("N8").Value = 91.038197
$ws.Range("O8").Value = 0.7437342022026642
$ws.Range("P8").Value = 0.744350442811843
$ws.Range("Q8").Value = 54.57315065230666
$ws.Range("R8").Value = 327.43890391384
$ws.Range("S8").Value = 0.007602669302746845
$ws.Range("T8").Value = 0.005089989585905238
$ws.Range("D9").Value = "MuSCs"
$ws.Range("G9").Value = 1.79836
$ws.Range("H9").Value = 3.59672
$ws.Range("I9").Value = 0.01022229350247785
$ws.Range("J9").Value = 0.006838162904394061
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.5
$ws.Range("M9").Value = 0.1013395
$ws.Range("N9").Value = 0.202679
$ws.Range("O9").Value = 0.002483671294065179
$ws.Range("P9").Value = 0.001657152803659562
$ws.Range("Q9").Value = 0.18224490322
$ws.Range("R9").Value = 0.72897961288
$ws.Range("S9").Value = [double]"2.538881693161324E-05"
$ws.Range("T9").Value = [double]"1.133188082889743E-05"
$ws.Range("G10").Value = 1.79836
$ws.Range("H10").Value = 3.59672
$ws.Range("I10").Value = 0.01022229350247785
$ws.Range("J10").Value = 0.006838162904394061
$ws.Range("M10").Value = 10.35489433333333
$ws.Range("N10").Value = 31.064683
$ws.Range("O10").Value = 0.2537821265032705
$ws.Range("P10").Value = 0.2539924043844974
$ws.Range("Q10").Value = 18.62182777329333
$ws.Range("R10").Value = 111.73096663976
$ws.Range("S10").Value = 0.002594235382799395
$ws.Range("T10").Value = 0.001736841437659925
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("G11").Value = 71.02944933333333
$ws.Range("H11").Value = 213.088348
$ws.Range("I11").Value = 0.4037477915460271
$ws.Range("J11").Value = 0.4051282381314676
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 30.34606566666666
$ws.Range("N11").Value = 91.038197
$ws.Range("O11").Value = 0.7437342022026642
$ws.Range("P11").Value = 0.744350442811843
$ws.Range("Q11").Value = 2155.464333736506
$ws.Range("R11").Value = 19399.17900362856
$ws.Range("S11").Value = 0.300281041636572
$ws.Range("T11").Value = 0.3015573834487397
$ws.Range("D12").Value = "MuSCs"
$ws.Range("G12").Value = 71.02944933333333
$ws.Range("H12").Value = 213.088348
$ws.Range("I12").Value = 0.4037477915460271
$ws.Range("J12").Value = 0.4051282381314676
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.5
$ws.Range("M12").Value = 0.1013395
$ws.Range("N12").Value = 0.202679
$ws.Range("O12").Value = 0.002483671294065179
$ws.Range("P12").Value = 0.001657152803659562
$ws.Range("Q12").Value = 7.198088880715333
$ws.Range("R12").Value = 43.188533284292
$ws.Range("S12").Value = 0.001002776799905079
$ws.Range("T12").Value = 0.0006713593956612203
$ws.Range("G13").Value = 71.02944933333333
$ws.Range("H13").Value = 213.088348
$ws.Range("I13").Value = 0.4037477915460271
$ws.Range("J13").Value = 0.4051282381314676
$ws.Range("M13").Value = 10.35489433333333
$ws.Range("N13").Value = 31.064683
$ws.Range("O13").Value = 0.2537821265032705
$ws.Range("P13").Value = 0.2539924043844974
$ws.Range("Q13").Value = 735.5024424015205
$ws.Range("R13").Value = 6619.521981613684
$ws.Range("S13").Value = 0.10246397310955
$ws.Range("T13").Value = 0.1028994952870667
